# Auto-generated edit script: updates F (想去人数 / interest count) and
# G (最低票价 / lowest price, row 9 on sheet 4 only) cell values to match
# the refreshed data snapshot described by the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 302
$ws.Cells.Item(5, 6).Value = 945
$ws.Cells.Item(7, 6).Value = 1535
$ws.Cells.Item(8, 6).Value = 39633
$ws.Cells.Item(8, 7).Value = 85
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(11, 6).Value = 8425
$ws.Cells.Item(15, 6).Value = 565
$ws.Cells.Item(18, 6).Value = 675
$ws.Cells.Item(19, 6).Value = 53
$ws.Cells.Item(21, 6).Value = 567
$ws.Cells.Item(22, 6).Value = 229
$ws.Cells.Item(23, 6).Value = 1072
$ws.Cells.Item(24, 6).Value = 349
$ws.Cells.Item(26, 6).Value = 389
$ws.Cells.Item(27, 6).Value = 567
$ws.Cells.Item(28, 6).Value = 585
$ws.Cells.Item(32, 6).Value = 5
$ws.Cells.Item(33, 6).Value = 395
$ws.Cells.Item(35, 6).Value = 141
$ws.Cells.Item(36, 6).Value = 833
$ws.Cells.Item(37, 6).Value = 372
$ws.Cells.Item(38, 6).Value = 16
$ws.Cells.Item(39, 6).Value = 176
$ws.Cells.Item(40, 6).Value = 51
$ws.Cells.Item(42, 6).Value = 1015
$ws.Cells.Item(43, 6).Value = 213
$ws.Cells.Item(44, 6).Value = 1041
$ws.Cells.Item(47, 6).Value = 9

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 4
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(6, 6).Value = 4384
$ws.Cells.Item(7, 6).Value = 7
$ws.Cells.Item(8, 6).Value = 302
$ws.Cells.Item(11, 6).Value = 71
$ws.Cells.Item(12, 6).Value = 92
$ws.Cells.Item(14, 6).Value = 46
$ws.Cells.Item(18, 6).Value = 160
$ws.Cells.Item(20, 6).Value = 4362

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 1787
$ws.Cells.Item(3, 6).Value = 427
$ws.Cells.Item(4, 6).Value = 332
$ws.Cells.Item(5, 6).Value = 102

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 1787
$ws.Cells.Item(3, 6).Value = 427
$ws.Cells.Item(4, 6).Value = 332
$ws.Cells.Item(5, 6).Value = 302
$ws.Cells.Item(7, 6).Value = 945
$ws.Cells.Item(8, 6).Value = 1535
$ws.Cells.Item(9, 6).Value = 39633
$ws.Cells.Item(9, 7).Value = 85
$ws.Cells.Item(11, 6).Value = 325
$ws.Cells.Item(13, 6).Value = 302
$ws.Cells.Item(14, 6).Value = 7
$ws.Cells.Item(15, 6).Value = 8426
$ws.Cells.Item(16, 6).Value = 158
$ws.Cells.Item(17, 6).Value = 540
$ws.Cells.Item(18, 6).Value = 71
$ws.Cells.Item(20, 6).Value = 565
$ws.Cells.Item(21, 6).Value = 92
$ws.Cells.Item(23, 6).Value = 212
$ws.Cells.Item(24, 6).Value = 675
$ws.Cells.Item(26, 6).Value = 53
$ws.Cells.Item(28, 6).Value = 567
$ws.Cells.Item(29, 6).Value = 229
$ws.Cells.Item(30, 6).Value = 1072
$ws.Cells.Item(32, 6).Value = 389
$ws.Cells.Item(33, 6).Value = 567
$ws.Cells.Item(34, 6).Value = 585
$ws.Cells.Item(37, 6).Value = 5
$ws.Cells.Item(38, 6).Value = 395
$ws.Cells.Item(39, 6).Value = 141
$ws.Cells.Item(40, 6).Value = 833
$ws.Cells.Item(41, 6).Value = 372
$ws.Cells.Item(42, 6).Value = 176
$ws.Cells.Item(44, 6).Value = 213
$ws.Cells.Item(45, 6).Value = 1041
$ws.Cells.Item(47, 6).Value = 74
